$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Create the new "ReportInfo" paragraph style (based on BodyTextIndent)
# ---------------------------------------------------------------------------
$reportInfo = $d.Styles.Add("ReportInfo", 1)
$reportInfo.BaseStyle = $d.Styles("BodyTextIndent")
$reportInfo.QuickStyle = $true
$reportInfo.ParagraphFormat.Alignment = 0
$reportInfo.Font.Name = "Courier New"
$reportInfo.Font.Size = 10

# ---------------------------------------------------------------------------
# 2. Insert three new paragraphs right after the "Generation 3" divider:
#      - empty "MainPersonText" paragraph
#      - "ReportInfo" paragraph with placeholder text
#      - empty "BodyTextIndent" paragraph
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Generation 3") {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Style = "MainPersonText"

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Style = "ReportInfo"
$p2.Range.Text = "This is used for Report Info section. "

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Style = "BodyTextIndent"

Write-Host "Done"
